$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new columns P and Q to the header row (14 and 15), copying style from O1 ---
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- New P/Q columns for data rows 2-25 are all zero (no special style) ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0
}

# --- Updated simulated loading-percent results for columns C(3) .. M(13), rows 2-25 ---
$newValues = @(
    @(3.722327560135659, 7.229550276666166, 11.62402334720648, 26.85825842887617, 31.59113135093992, 1.764284764925966, 2.802000026140288, 10.90294933493978, 18.30768812800111, 7.195790056695235, 24.21210140620341),
    @(3.436742906384903, 6.846544035320609, 11.02352491914789, 26.21113797823519, 30.80265726062291, 2.050103161057943, 2.564043389562335, 10.87011219372632, 17.99905784558912, 6.95412115394206, 22.61517086919804),
    @(3.247782839637662, 6.602064312215719, 10.63894869504672, 25.81620704573505, 30.32656604665539, 2.231425215488349, 2.532203910772701, 10.85432496729646, 17.81200806643223, 6.800855269192743, 21.57615535257406),
    @(3.154579986502397, 6.504032153730228, 10.47978380038966, 25.64067276221575, 30.10728407495156, 2.307207073216465, 2.601178426204153, 10.84402121369418, 17.7226075655863, 6.736658193660608, 21.1377815501557),
    @(3.123355595626776, 6.492362357451513, 10.45491549354487, 25.59283107497871, 30.03750168895643, 2.320363717598884, 2.617140197009079, 10.83635484522943, 17.69091045357851, 6.725257327229293, 21.06385737246452),
    @(3.205262701117521, 6.613455182079095, 10.64163058363756, 25.76305286827741, 30.2327916583165, 2.233741663142332, 2.544944044301866, 10.83777565983411, 17.76488148576699, 6.798189909713827, 21.56965955222175),
    @(3.576834634318336, 7.115183444328581, 11.42629086692671, 26.56947917254366, 31.2010082372476, 1.863298092810945, 2.72229120163278, 10.86914737755857, 18.14117379012997, 7.111221487778375, 23.67276140083864),
    @(4.264752864991602, 8.002562833597924, 12.82875868152507, 28.25211878103336, 33.31511303636957, 1.932504511756898, 3.292344901740542, 10.99663621091517, 18.97761597358864, 7.688069539424195, 27.33704760767613),
    @(4.730714768107982, 8.54135804975564, 13.56224729089071, 29.18885937395603, 34.48193041080436, 2.382654150638599, 3.686815296516222, 11.03412949355121, 19.43193176051998, 7.972459553229297, 29.74906234742395),
    @(5.201365732084928, 8.190442939291307, 12.04611256857507, 27.14391975519309, 31.65491122766553, 3.176031790840131, 3.805316413473979, 10.39795920319963, 18.29232772845502, 7.176183847190372, 30.75746828973589),
    @(5.561500931913464, 7.763597007811438, 10.65354509364991, 25.27203086291273, 29.11114145800326, 4.352325631379124, 3.826741297714936, 9.878534596766315, 17.29222693586744, 6.532023791051563, 31.12062152406765),
    @(5.824819100073502, 7.25632276591338, 9.269665360402476, 23.32551385647625, 26.46012888740123, 5.61791408178646, 3.776105208381, 9.394048586349053, 16.26322906772045, 5.982518611299803, 31.02072496045494),
    @(5.961718066192073, 6.862170371205227, 8.32863507557588, 21.9168987822018, 24.52673537487875, 6.536425104918109, 3.708994508072371, 9.072160033361854, 15.52010317421031, 5.673220207895515, 30.7402453230033),
    @(5.96164532937109, 6.751781217701046, 8.096210232483122, 21.54498169116122, 24.00812376519009, 6.744319923308831, 3.677250194686278, 8.997080905398786, 15.31973606983043, 5.606491974635159, 30.57140522994397),
    @(5.725022326415751, 6.602777400548366, 7.984781068159395, 21.45641326191044, 23.89476277394383, 6.460031588123118, 3.52086943580349, 9.072151825661397, 15.28682727151859, 5.576282893678498, 29.60247934117938),
    @(5.459047071780731, 6.706160774082731, 8.399306841222028, 22.14454416558284, 24.85312958111439, 5.701471837324598, 3.436246845906917, 9.297265868627397, 15.65797300058517, 5.702647999561766, 28.99862680040538),
    @(5.171167551565143, 7.033107237658526, 9.3915339789291, 23.6216884451322, 26.90517900218183, 4.483610673149109, 3.404561842668864, 9.69423069086737, 16.45536540115175, 6.062834282806045, 28.65537762165196),
    @(4.837017594967094, 7.531362117595515, 10.84864978375311, 25.55540360057255, 29.54162181966277, 3.13690534721214, 3.426519742205345, 10.18632185035538, 17.47056403718411, 6.680308967359077, 28.55263741611301),
    @(4.51926754922526, 8.431131921635478, 13.37558562036826, 28.78860354331195, 33.90416645649819, 2.260555748785416, 3.588153002665651, 10.96989389123893, 19.17303325548318, 7.889138552941042, 29.13213488421318),
    @(4.828337960117391, 8.930127375129677, 14.24985391813835, 29.97428756173824, 35.45461489970504, 2.624388140447137, 3.895768085193564, 11.13292032574677, 19.7868438071571, 8.275627846155928, 30.9029699972947),
    @(5.074377143087956, 9.19892380238143, 14.68890435052573, 30.66023384594697, 36.37666142074222, 2.851259745528094, 4.089224689637798, 11.22751954475032, 20.16581145141496, 8.46706753007261, 32.01052682309825),
    @(4.977274707231053, 9.044943870677104, 14.45158727866099, 30.34567184450714, 35.97545900436108, 2.731793435761431, 3.984783912374489, 11.19496340188318, 20.01204871850687, 8.367080241627505, 31.42451648703549),
    @(4.539593378994615, 8.460910485105607, 13.52961419717534, 29.07382171589499, 34.325521074649, 2.269872137517071, 3.585803822691178, 11.05448257897538, 19.35836379785382, 7.975869970076581, 29.10624828290923),
    @(4.029197081185658, 7.790538984672701, 12.47089876473525, 27.71095802328734, 32.58542060186127, 1.754330402511626, 3.143557848181273, 10.92884558282484, 18.67060005471902, 7.533549306613914, 26.39356237998738)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $newValues[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 3
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
